$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the edited cells keep a text (string) format so that values such
# as "2345" are not reinterpreted as numbers.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("F4").NumberFormat = "@"

# Row 3 (Id = 2)
$ws.Range("B3").Value = "Nguyễn Khôi"
$ws.Range("C3").Value = "2345"
$ws.Range("D3").Value = "Tỉnh Cao Bằng-Huyện Bảo Lâm-Xã Vĩnh Quang"
$ws.Range("F3").Value = "245236587"

# Row 4 (Id = 3)
$ws.Range("B4").Value = "Trần Khánh"
$ws.Range("D4").Value = "Tỉnh Hà Giang-Huyện Vị Xuyên-Thị trấn Nông Trường Việt Lâm"
$ws.Range("F4").Value = "123321658"
